# The "Commands" sheet listed several Artisan commands that have since been
# removed from the application: keyboard(<bool>), showCurve(<name>,<bool>),
# showExtraCurve(<extra_device>,<curve>,<bool>), showEvents(<event_type>,<bool>)
# and showBackgroundEvents(<bool>). These occupied rows 95-99 (5 rows) just
# before the "RC Command" section. Remove them; everything below shifts up.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Rows("95:99").Delete()

# Update the view state to reflect the new layout (top of sheet, selection on A3).
$ws.Activate()
$ws.Range("A3").Select()
